$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (avoids Excel auto-converting
# numeric-looking strings like "7.33" or "0.0931" into real numbers),
# then restores the cell to its original (unstyled) Normal style so no
# stray formatting is left behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "63.607.57"
Set-TextValue $ws.Range("E2") "  +2.59%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.477.56"
Set-TextValue $ws.Range("E3") "  +2.39%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.20%  "

# Row 5
Set-TextValue $ws.Range("D5") "576.38"
Set-TextValue $ws.Range("E5") "  +2.38%  "

# Row 6
Set-TextValue $ws.Range("D6") "149.12"
Set-TextValue $ws.Range("E6") "  +4.05%  "

# Row 7
Set-TextValue $ws.Range("E7") "  -0.08%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.541"
Set-TextValue $ws.Range("E8") "  +1.80%  "

# Row 9
Set-TextValue $ws.Range("E9") "  +4.77%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.155"
Set-TextValue $ws.Range("E10") "  +0.77%  "

# Row 11
Set-TextValue $ws.Range("B11") "Cardano"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D11") "0.364"
Set-TextValue $ws.Range("E11") "  +4.18%  "

# Row 12
Set-TextValue $ws.Range("B12") "Toncoin"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D12") "5.35"
Set-TextValue $ws.Range("E12") "  +3.07%  "

# Row 13
Set-TextValue $ws.Range("D13") "27.39"
Set-TextValue $ws.Range("E13") "  +5.03%  "

# Row 14
Set-TextValue $ws.Range("D14") "0.0000186"
Set-TextValue $ws.Range("E14") "  +7.18%  "

# Row 15
Set-TextValue $ws.Range("D15") "2.948.97"
Set-TextValue $ws.Range("E15") "  +3.31%  "

# Row 16
Set-TextValue $ws.Range("D16") "63.578.19"
Set-TextValue $ws.Range("E16") "  +2.75%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.486.26"
Set-TextValue $ws.Range("E17") "  +3.54%  "

# Row 18
Set-TextValue $ws.Range("D18") "11.59"
Set-TextValue $ws.Range("E18") "  +2.01%  "

# Row 19
Set-TextValue $ws.Range("D19") "7.33"
Set-TextValue $ws.Range("E19") "  +7.32%  "

# Row 20
Set-TextValue $ws.Range("D20") "4.25"
Set-TextValue $ws.Range("E20") "  +2.92%  "

# Row 21
Set-TextValue $ws.Range("D21") "329.43"
Set-TextValue $ws.Range("E21") "  +1.93%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.999"
Set-TextValue $ws.Range("E22") "  -0.03%  "

# Row 23
Set-TextValue $ws.Range("D23") "1.92"
Set-TextValue $ws.Range("E23") "  +10.66%  "

# Row 24
Set-TextValue $ws.Range("D24") "67.57"
Set-TextValue $ws.Range("E24") "  +1.13%  "

# Row 25
Set-TextValue $ws.Range("D25") "632.29"
Set-TextValue $ws.Range("E25") "  +13.82%  "

# Row 26
Set-TextValue $ws.Range("E26") "  +13.30%  "

# Row 27
Set-TextValue $ws.Range("D27") "8.79"
Set-TextValue $ws.Range("E27") "  +0.41%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.605.15"
Set-TextValue $ws.Range("E28") "  +2.61%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.53"
Set-TextValue $ws.Range("E29") "  +9.80%  "

# Row 30
Set-TextValue $ws.Range("D30") "8.49"
Set-TextValue $ws.Range("E30") "  +3.42%  "

# Row 31
Set-TextValue $ws.Range("E31") "  -0.23%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.144"
Set-TextValue $ws.Range("E32") "  -1.88%  "

# Row 33
Set-TextValue $ws.Range("E33") "  +2.40%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.22"
Set-TextValue $ws.Range("E34") "  +10.01%  "

# Row 35
Set-TextValue $ws.Range("E35") "  +3.81%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.998"
Set-TextValue $ws.Range("E36") "  -0.14%  "

# Row 37
Set-TextValue $ws.Range("E37") "  +2.17%  "

# Row 38
Set-TextValue $ws.Range("D38") "5.56"
Set-TextValue $ws.Range("E38") "  +2.54%  "

# Row 39
Set-TextValue $ws.Range("D39") "19.01"
Set-TextValue $ws.Range("E39") "  +2.54%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.86"
Set-TextValue $ws.Range("E40") "  +2.34%  "

# Row 41
Set-TextValue $ws.Range("D41") "147.47"
Set-TextValue $ws.Range("E41") "  -3.91%  "

# Row 42
Set-TextValue $ws.Range("E42") "  +19.15%  "

# Row 43
Set-TextValue $ws.Range("E43") "  -0.02%  "

# Row 44
Set-TextValue $ws.Range("D44") "151.35"
Set-TextValue $ws.Range("E44") "  +2.86%  "

# Row 45
Set-TextValue $ws.Range("D45") "3.77"
Set-TextValue $ws.Range("E45") "  +3.66%  "

# Row 46
Set-TextValue $ws.Range("D46") "21.26"
Set-TextValue $ws.Range("E46") "  +7.40%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.0551"
Set-TextValue $ws.Range("E47") "  +4.42%  "

# Row 48
Set-TextValue $ws.Range("E48") "  +3.20%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0240"
Set-TextValue $ws.Range("E49") "  +5.50%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.0931"
Set-TextValue $ws.Range("E50") "  +1.07%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.750"
Set-TextValue $ws.Range("E51") "  +5.10%  "

